$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report date range) ---
$ws.Range("A8").Value = "Volume 32   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/1/2025  Through  12/7/2025"

# --- Weekly crime statistics table (rows 14-33) ---
# Row 14
$ws.Range("F14").Value = 2
$ws.Range("I14").Value = 36
$ws.Range("K14").Value = -21.739130434782
$ws.Range("L14").Value = -38.983050847457
$ws.Range("M14").Value = -56.097560975609
$ws.Range("N14").Value = -84.680851063829

# Row 15
$ws.Range("C15").Value = 5
$ws.Range("E15").Value = 400
$ws.Range("F15").Value = 19
$ws.Range("G15").Value = 15
$ws.Range("H15").Value = 26.666666666666
$ws.Range("I15").Value = 247
$ws.Range("J15").Value = 214
$ws.Range("K15").Value = 15.420560747663
$ws.Range("L15").Value = 18.75
$ws.Range("M15").Value = 47.904191616766
$ws.Range("N15").Value = -52.681992337164

# Row 16
$ws.Range("C16").Value = 27
$ws.Range("E16").Value = -22.857142857142
$ws.Range("F16").Value = 129
$ws.Range("G16").Value = 122
$ws.Range("H16").Value = 5.737704918032
$ws.Range("I16").Value = 1558
$ws.Range("J16").Value = 1619
$ws.Range("K16").Value = -3.767757875231
$ws.Range("L16").Value = -6.144578313253
$ws.Range("M16").Value = -42.972181551976
$ws.Range("N16").Value = -88.129523809523

# Row 17
$ws.Range("C17").Value = 65
$ws.Range("D17").Value = 67
$ws.Range("E17").Value = -2.985074626865
$ws.Range("F17").Value = 290
$ws.Range("G17").Value = 247
$ws.Range("H17").Value = 17.408906882591
$ws.Range("I17").Value = 3710
$ws.Range("J17").Value = 3480
$ws.Range("K17").Value = 6.609195402298
$ws.Range("L17").Value = 11.578947368421
$ws.Range("M17").Value = 59.638554216867
$ws.Range("N17").Value = -41.501103752759

# Row 18
$ws.Range("C18").Value = 26
$ws.Range("D18").Value = 23
$ws.Range("E18").Value = 13.043478260869
$ws.Range("F18").Value = 113
$ws.Range("G18").Value = 117
$ws.Range("H18").Value = -3.418803418803
$ws.Range("I18").Value = 1375
$ws.Range("J18").Value = 1447
$ws.Range("K18").Value = -4.975812024879
$ws.Range("L18").Value = -13.249211356466
$ws.Range("M18").Value = -57.796193984039
$ws.Range("N18").Value = -92.108585858585

# Row 19
$ws.Range("C19").Value = 110
$ws.Range("D19").Value = 107
$ws.Range("E19").Value = 2.803738317757
$ws.Range("F19").Value = 415
$ws.Range("G19").Value = 382
$ws.Range("H19").Value = 8.638743455497
$ws.Range("I19").Value = 5324
$ws.Range("J19").Value = 5438
$ws.Range("K19").Value = -2.096358955498
$ws.Range("L19").Value = -13.006535947712
$ws.Range("M19").Value = 2.740254727904
$ws.Range("N19").Value = -36.588851834206

# Row 20
$ws.Range("C20").Value = 44
$ws.Range("D20").Value = 39
$ws.Range("E20").Value = 12.820512820512
$ws.Range("F20").Value = 138
$ws.Range("G20").Value = 165
$ws.Range("H20").Value = -16.363636363636
$ws.Range("I20").Value = 1712
$ws.Range("J20").Value = 1924
$ws.Range("K20").Value = -11.018711018711
$ws.Range("L20").Value = -1.89111747851
$ws.Range("M20").Value = -4.941699056079
$ws.Range("N20").Value = -92.155065756312

# Row 21
$ws.Range("C21").Value = 277
$ws.Range("D21").Value = 272
$ws.Range("E21").Value = 1.838235294117
$ws.Range("F21").Value = 1106
$ws.Range("G21").Value = 1048
$ws.Range("H21").Value = 5.534351145038
$ws.Range("I21").Value = 13962
$ws.Range("J21").Value = 14168
$ws.Range("K21").Value = -1.453980801806
$ws.Range("L21").Value = -5.033328798802
$ws.Range("M21").Value = -10.189116171362
$ws.Range("N21").Value = -79.427409492094

# Row 22
$ws.Range("C22").Value = 6
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 13
$ws.Range("H22").Value = -35
$ws.Range("I22").Value = 216
$ws.Range("J22").Value = 194
$ws.Range("K22").Value = 11.340206185567
$ws.Range("L22").Value = 11.917098445595
$ws.Range("M22").Value = -16.602316602316

# Row 23
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 10
$ws.Range("E23").Value = -10
$ws.Range("F23").Value = 40
$ws.Range("G23").Value = 38
$ws.Range("H23").Value = 5.263157894736
$ws.Range("I23").Value = 456
$ws.Range("J23").Value = 473
$ws.Range("K23").Value = -3.594080338266
$ws.Range("L23").Value = -7.505070993914
$ws.Range("M23").Value = 42.5

# Row 24
$ws.Range("C24").Value = 275
$ws.Range("D24").Value = 287
$ws.Range("E24").Value = -4.181184668989
$ws.Range("F24").Value = 1032
$ws.Range("G24").Value = 1141
$ws.Range("H24").Value = -9.553023663453
$ws.Range("I24").Value = 12615
$ws.Range("J24").Value = 14017
$ws.Range("K24").Value = -10.002140258257
$ws.Range("L24").Value = -14.445574771108
$ws.Range("M24").Value = 8.144020574367

# Row 25
$ws.Range("C25").Value = 99
$ws.Range("D25").Value = 121
$ws.Range("E25").Value = -18.181818181818
$ws.Range("F25").Value = 353
$ws.Range("G25").Value = 521
$ws.Range("H25").Value = -32.245681381957
$ws.Range("I25").Value = 5017
$ws.Range("J25").Value = 6658
$ws.Range("K25").Value = -24.647041153499
$ws.Range("L25").Value = -29.546412020783

# Row 26
$ws.Range("C26").Value = 113
$ws.Range("D26").Value = 131
$ws.Range("E26").Value = -13.740458015267
$ws.Range("F26").Value = 490
$ws.Range("G26").Value = 511
$ws.Range("H26").Value = -4.109589041095
$ws.Range("I26").Value = 6050
$ws.Range("J26").Value = 6155
$ws.Range("K26").Value = -1.705930138099
$ws.Range("L26").Value = 9.324177809902
$ws.Range("M26").Value = -4.044409199048

# Row 27
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 500
$ws.Range("F27").Value = 21
$ws.Range("G27").Value = 16
$ws.Range("H27").Value = 31.25
$ws.Range("I27").Value = 286
$ws.Range("J27").Value = 303
$ws.Range("K27").Value = -5.610561056105
$ws.Range("L27").Value = -5.921052631578

# Row 28
$ws.Range("C28").Value = 10
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = 66.666666666666
$ws.Range("F28").Value = 55
$ws.Range("G28").Value = 42
$ws.Range("H28").Value = 30.952380952381
$ws.Range("I28").Value = 633
$ws.Range("J28").Value = 636
$ws.Range("K28").Value = -0.471698113207
$ws.Range("L28").Value = 5.852842809364

# Row 29
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("D29").Value = 6
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E29").Value = -66.666666666666
$ws.Range("F29").Value = 7
$ws.Range("G29").Value = 11
$ws.Range("H29").Value = -36.363636363636
$ws.Range("I29").Value = 123
$ws.Range("J29").Value = 103
$ws.Range("K29").Value = 19.417475728155
$ws.Range("L29").Value = -15.172413793103
$ws.Range("M29").Value = -53.231939163498
$ws.Range("N29").Value = -83.333333333333

# Row 30
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 5
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E30").Value = -60
$ws.Range("F30").Value = 7
$ws.Range("G30").Value = 9
$ws.Range("H30").Value = -22.222222222222
$ws.Range("I30").Value = 96
$ws.Range("J30").Value = 88
$ws.Range("K30").Value = 9.090909090909
$ws.Range("L30").Value = -21.951219512195
$ws.Range("M30").Value = -55.760368663594
$ws.Range("N30").Value = -85.093167701863

# Row 31
$ws.Range("D31").Value = 1
$ws.Range("F31").Value = 3
$ws.Range("G31").Value = 6
$ws.Range("H31").Value = -50
$ws.Range("J31").Value = 145
$ws.Range("K31").Value = -25.51724137931
$ws.Range("L31").Value = -1.818181818181

# Row 33
$ws.Range("D33").Value = 2
$ws.Range("G33").Value = 6
$ws.Range("J33").Value = 42
$ws.Range("K33").Value = -11.904761904761

